$d = $word.ActiveDocument

# Locate the paragraph that contains the "Requisitos" marker text so we
# find the right block even if paragraph indices shift for any reason.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Ver no Jupiter*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $startRange = $target.Range
    $start = $startRange.Start

    # The paragraph after "Ver no Jupiter..." holds the copyright notice;
    # the one after that is a blank separator paragraph. Remove all three
    # (the "Ver no Jupiter..." paragraph, the "c 2020 ..." paragraph, and
    # the trailing blank paragraph) while keeping the blank paragraph that
    # precedes them and the page-break paragraph that follows.
    $copyrightPara = $target.Next()
    $blankPara = $copyrightPara.Next()
    $end = $blankPara.Range.End

    $deleteRange = $d.Range($start, $end)
    $deleteRange.Delete()
}
